# Update the "Förändrad" (Changed) date column C for rows 2-23
# from 2023-09-15 (45184) to 2023-09-16 (45185), i.e. +1 day.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 23; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45184) {
        $cell.Value2 = 45185
    }
}
